$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.189.30"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "3.372.47"
$ws.Range("E3").Value = "  +1.61%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D5").Value = "'572.45"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").Value = "'137.40"
$ws.Range("E6").Value = "  +7.57%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "3.373.21"
$ws.Range("E8").Value = "  +1.70%  "
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("D10").Value = "'7.61"
$ws.Range("E10").Value = "  +3.62%  "
$ws.Range("E11").Value = "  +5.03%  "
$ws.Range("E12").Value = "  +5.41%  "
$ws.Range("D13").Value = "3.951.47"
$ws.Range("E13").Value = "  +1.44%  "
$ws.Range("E14").Value = "  +1.92%  "
$ws.Range("E15").Value = "  +3.34%  "
$ws.Range("D16").Value = "3.371.64"
$ws.Range("E16").Value = "  +1.27%  "
$ws.Range("D17").Value = "'25.17"
$ws.Range("E17").Value = "  +1.63%  "
$ws.Range("D18").Value = "61.218.18"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").Value = "'13.91"
$ws.Range("E19").Value = "  +4.77%  "
$ws.Range("E20").Value = "  +4.73%  "
$ws.Range("D21").Value = "'9.40"
$ws.Range("E21").Value = "  +5.07%  "
$ws.Range("D22").Value = "'379.24"
$ws.Range("E22").Value = "  +7.17%  "
$ws.Range("E23").Value = "  +3.07%  "
$ws.Range("D24").Value = "3.509.11"
$ws.Range("E24").Value = "  +1.60%  "
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").Value = "'70.70"
$ws.Range("E26").Value = "  +1.72%  "
$ws.Range("D27").Value = "'0.0000119"
$ws.Range("E27").Value = "  +11.26%  "
$ws.Range("D28").Value = "'1.66"
$ws.Range("E28").Value = "  +17.55%  "
$ws.Range("E29").Value = "  +8.52%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").Value = "  +4.19%  "
$ws.Range("E32").Value = "  +5.05%  "
$ws.Range("E33").Value = "  +0.86%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").Value = "3.404.53"
$ws.Range("E35").Value = "  +1.53%  "
$ws.Range("D36").Value = "'23.38"
$ws.Range("E36").Value = "  +3.97%  "
$ws.Range("E37").Value = "  +7.40%  "
$ws.Range("E38").Value = "  +4.59%  "
$ws.Range("E39").Value = "  +5.07%  "
$ws.Range("D40").Value = "'162.33"
$ws.Range("E40").Value = "  +0.43%  "
$ws.Range("D41").Value = "'0.0796"
$ws.Range("E41").Value = "  +5.49%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").Value = "'1.71"
$ws.Range("E43").Value = "  +9.54%  "
$ws.Range("D44").Value = "'4.43"
$ws.Range("E44").Value = "  +2.89%  "
$ws.Range("D45").Value = "'41.52"
$ws.Range("E45").Value = "  +0.71%  "
$ws.Range("E46").Value = "  +7.29%  "
$ws.Range("E47").Value = "  +1.86%  "
$ws.Range("D48").Value = "'23.35"
$ws.Range("E48").Value = "  +5.44%  "
$ws.Range("E49").Value = "  +4.50%  "
$ws.Range("D50").Value = "'23.01"
$ws.Range("E50").Value = "  +9.04%  "
$ws.Range("D51").Value = "2.344.03"
$ws.Range("E51").Value = "  +7.63%  "
